# Apply weekly-report update:
#  - New "Friday (07/11/2025)" day-block inserted before the existing
#    Saturday block (pushes Saturday/Sunday blocks down by 6 rows).
#  - Report-generated timestamp, billed total and line-item count updated.
#  - All "Pricing" ($) values across every day block zeroed out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the new Friday block: insert 6 blank rows at row 14.
#    Everything that was at rows 14-27 (Saturday + Sunday blocks) shifts
#    down to rows 20-33; merged-cell ranges shift automatically too.
# ---------------------------------------------------------------------
$ws.Rows("14:19").Insert()

# ---------------------------------------------------------------------
# 2) Re-create the Friday block in the freshly emptied rows 14-17 by
#    copying formatting from the (now shifted) Saturday block, so the
#    banding / fonts / fills match exactly, then overwrite the text.
# ---------------------------------------------------------------------
$ws.Range("A20:I20").Copy($ws.Range("A14"))   # day-header row style
$ws.Range("A21:I21").Copy($ws.Range("A15"))   # column-header row
$ws.Range("A22:I22").Copy($ws.Range("A16"))   # first point-row style
$ws.Range("A26:I26").Copy($ws.Range("A17"))   # TOTAL row style

$ws.Range("A14").Value = "Friday (07/11/2025)"

$ws.Range("A16").Value = "Point 08"
$ws.Range("B16").Value = "PLA-HDIG"
$ws.Range("C16").Value = "Inst"
$ws.Range("D16").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("E16").Value = "EA"
$ws.Range("F16").Value = 0
# G16 is left as-is: the copy from A22:I22 already carried over the
# blank-inline-string "N/A" placeholder cell used by every point row.
$ws.Range("H16").Value = 0

$ws.Range("A17").Value = "TOTAL"
$ws.Range("H17").Value = 0

$ws.Range("A17:G17").Merge()

# ---------------------------------------------------------------------
# 3) Zero out every "Pricing" figure in the shifted Saturday & Sunday
#    blocks (the units/#s stay as they were).
# ---------------------------------------------------------------------
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0

$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0

# ---------------------------------------------------------------------
# 4) Header / summary field updates.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 7
